$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Runs" (sheet1.xml)
# Row 1 (Run105) values are updated in place.
# Rows 2-5 are new rows (Run108, Run109, Run110, Run111) inserted before
# the old row 2 (Run116), which becomes row 6 with its original values.
# ----------------------------------------------------------------------
$wsRuns = $wb.Worksheets.Item("Runs")

$runsData = @(
    @("Run105", "0.560 (0.449)", "0.455 (0.287)", "0.600 (0.491)", "0.560 (0.428)", "0.580 (0.244)", "0.550 (0.499)", "0.768 (0.386)", "0.716 (0.210)", "0.779 (0.415)", "0.636 (0.435)", "0.555 (0.280)", "0.664 (0.473)"),
    @("Run108", "0.509 (0.444)", "0.758 (0.160)", "0.382 (0.486)", "0.617 (0.456)", "0.793 (0.197)", "0.564 (0.497)", "0.759 (0.414)", "0.607 (0.324)", "0.774 (0.419)", "0.613 (0.450)", "0.745 (0.202)", "0.572 (0.495)"),
    @("Run109", "0.743 (0.383)", "0.733 (0.173)", "0.746 (0.436)", "0.776 (0.390)", "0.885 (0.100)", "0.741 (0.440)", "0.901 (0.260)", "0.801 (0.087)", "0.918 (0.275)", "0.807 (0.352)", "0.779 (0.156)", "0.814 (0.389)"),
    @("Run110", "0.720 (0.439)", "0.722 (0.169)", "0.719 (0.450)", "0.864 (0.299)", "0.875 (0.132)", "0.859 (0.349)", "0.816 (0.341)", "0.803 (0.141)", "0.820 (0.384)", "0.781 (0.386)", "0.809 (0.150)", "0.775 (0.418)"),
    @("Run111", "0.775 (0.394)", "0.872 (0.074)", "0.752 (0.432)", "0.810 (0.341)", "0.826 (0.159)", "0.802 (0.399)", "0.879 (0.307)", "0.676 (0.220)", "0.895 (0.307)", "0.819 (0.358)", "0.826 (0.152)", "0.818 (0.386)"),
    @("Run116", "0.130 (0.214)", "0.399 (0.183)", "0.000 (0.000)", "0.250 (0.323)", "0.531 (0.216)", "0.052 (0.223)", "0.142 (0.288)", "0.549 (0.176)", "0.053 (0.224)", "0.151 (0.266)", "0.466 (0.201)", "0.030 (0.172)")
)

for ($i = 0; $i -lt $runsData.Length; $i++) {
    $rowNum = $i + 1
    $rowVals = $runsData[$i]
    $wsRuns.Cells.Item($rowNum, 1).Value = 0
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $wsRuns.Cells.Item($rowNum, $c + 2).Value = $rowVals[$c]
    }
}

# ----------------------------------------------------------------------
# Sheet "GAN" (sheet3.xml)
# Previously empty, now contains one data row (Iteration30).
# ----------------------------------------------------------------------
$wsGan = $wb.Worksheets.Item("GAN")

$ganRow = @("Iteration30", "0.69 (0.10)", "0.77 (0.06)", "0.65 (0.16)", "0.77 (0.09)", "0.84 (0.04)", "0.74 (0.11)", "0.84 (0.06)", "0.72 (0.08)", "0.85 (0.06)", "0.75 (0.08)", "0.79 (0.03)", "0.74 (0.10)")

$wsGan.Cells.Item(1, 1).Value = 0
for ($c = 0; $c -lt $ganRow.Length; $c++) {
    $wsGan.Cells.Item(1, $c + 2).Value = $ganRow[$c]
}
